# TensorVals.xlsx — "Add files via upload": insert a new Sheet2 with a
# Bi / lambda1 / A1 lookup table (transient-conduction Heisler chart data),
# make it the active sheet, and leave Sheet1 selection parked at C21.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes right after Sheet1 (tab order Sheet1, Sheet2).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row, styled with the built-in "Note" cell style (pale-yellow fill
# + thin gray border).
$ws2.Cells.Item(1,1).Value = "Bi"
$ws2.Cells.Item(1,2).Value = "lambda1"
$ws2.Cells.Item(1,3).Value = "A1"
$ws2.Range("A1:C1").Style = "Note"

# Data rows.
$ws2.Cells.Item(2,1).Value = 0.01
$ws2.Cells.Item(2,2).Value = 0.14119999999999999
$ws2.Cells.Item(2,3).Value = 1.0024999999999999
$ws2.Cells.Item(3,1).Value = 0.02
$ws2.Cells.Item(3,2).Value = 0.19950000000000001
$ws2.Cells.Item(3,3).Value = 1.0049999999999999
$ws2.Cells.Item(4,1).Value = 0.04
$ws2.Cells.Item(4,2).Value = 0.28139999999999998
$ws2.Cells.Item(4,3).Value = 1.0099
$ws2.Cells.Item(5,1).Value = 0.06
$ws2.Cells.Item(5,2).Value = 0.34379999999999999
$ws2.Cells.Item(5,3).Value = 1.1048
$ws2.Cells.Item(6,1).Value = 0.08
$ws2.Cells.Item(6,2).Value = 0.39600000000000002
$ws2.Cells.Item(6,3).Value = 1.0197000000000001
$ws2.Cells.Item(7,1).Value = 0.1
$ws2.Cells.Item(7,2).Value = 0.44169999999999998
$ws2.Cells.Item(7,3).Value = 1.0246
$ws2.Cells.Item(8,1).Value = 0.2
$ws2.Cells.Item(8,2).Value = 0.61699999999999999
$ws2.Cells.Item(8,3).Value = 1.0483
$ws2.Cells.Item(9,1).Value = 0.3
$ws2.Cells.Item(9,2).Value = 0.74650000000000005
$ws2.Cells.Item(9,3).Value = 1.0711999999999999
$ws2.Cells.Item(10,1).Value = 0.4
$ws2.Cells.Item(10,2).Value = 0.85160000000000002
$ws2.Cells.Item(10,3).Value = 1.0931
$ws2.Cells.Item(11,1).Value = 0.5
$ws2.Cells.Item(11,2).Value = 0.94079999999999997
$ws2.Cells.Item(11,3).Value = 1.1143000000000001
$ws2.Cells.Item(12,1).Value = 0.6
$ws2.Cells.Item(12,2).Value = 1.0184
$ws2.Cells.Item(12,3).Value = 1.1345000000000001
$ws2.Cells.Item(13,1).Value = 0.7
$ws2.Cells.Item(13,2).Value = 1.0872999999999999
$ws2.Cells.Item(13,3).Value = 1.1538999999999999
$ws2.Cells.Item(14,1).Value = 0.8
$ws2.Cells.Item(14,2).Value = 1.149
$ws2.Cells.Item(14,3).Value = 1.1724000000000001
$ws2.Cells.Item(15,1).Value = 0.9
$ws2.Cells.Item(15,2).Value = 1.2048000000000001
$ws2.Cells.Item(15,3).Value = 1.1901999999999999
$ws2.Cells.Item(16,1).Value = 1
$ws2.Cells.Item(16,2).Value = 1.2558
$ws2.Cells.Item(16,3).Value = 1.2071000000000001
$ws2.Cells.Item(17,1).Value = 2
$ws2.Cells.Item(17,2).Value = 1.5994999999999999
$ws2.Cells.Item(17,3).Value = 1.3384
$ws2.Cells.Item(18,1).Value = 3
$ws2.Cells.Item(18,2).Value = 1.7887
$ws2.Cells.Item(18,3).Value = 1.4191
$ws2.Cells.Item(19,1).Value = 4
$ws2.Cells.Item(19,2).Value = 1.9080999999999999
$ws2.Cells.Item(19,3).Value = 1.4698
$ws2.Cells.Item(20,1).Value = 5
$ws2.Cells.Item(20,2).Value = 1.9898
$ws2.Cells.Item(20,3).Value = 1.5028999999999999
$ws2.Cells.Item(21,1).Value = 6
$ws2.Cells.Item(21,2).Value = 2.0489999999999999
$ws2.Cells.Item(21,3).Value = 1.5253000000000001
$ws2.Cells.Item(22,1).Value = 7
$ws2.Cells.Item(22,2).Value = 2.0937000000000001
$ws2.Cells.Item(22,3).Value = 1.5410999999999999
$ws2.Cells.Item(23,1).Value = 8
$ws2.Cells.Item(23,2).Value = 2.1286
$ws2.Cells.Item(23,3).Value = 1.5526
$ws2.Cells.Item(24,1).Value = 9
$ws2.Cells.Item(24,2).Value = 2.1566000000000001
$ws2.Cells.Item(24,3).Value = 1.5610999999999999
$ws2.Cells.Item(25,1).Value = 10
$ws2.Cells.Item(25,2).Value = 2.1795
$ws2.Cells.Item(25,3).Value = 1.5677000000000001
$ws2.Cells.Item(26,1).Value = 20
$ws2.Cells.Item(26,2).Value = 2.2879999999999998
$ws2.Cells.Item(26,3).Value = 1.5919000000000001
$ws2.Cells.Item(27,1).Value = 30
$ws2.Cells.Item(27,2).Value = 2.3260999999999998
$ws2.Cells.Item(27,3).Value = 1.5972999999999999
$ws2.Cells.Item(28,1).Value = 40
$ws2.Cells.Item(28,2).Value = 2.4550000000000001
$ws2.Cells.Item(28,3).Value = 1.5992999999999999
$ws2.Cells.Item(29,1).Value = 50
$ws2.Cells.Item(29,2).Value = 2.3572000000000002
$ws2.Cells.Item(29,3).Value = 1.6002099999999999
$ws2.Cells.Item(30,1).Value = 100
$ws2.Cells.Item(30,2).Value = 2.3809
$ws2.Cells.Item(30,3).Value = 1.6014999999999999
$ws2.Cells.Item(31,1).Value = 100000
$ws2.Cells.Item(31,2).Value = 2.4047999999999998
$ws2.Cells.Item(31,3).Value = 1.6021000000000001

# Column widths, matching the auto-fit sizing seen in the saved workbook.
$ws2.Columns.Item(1).ColumnWidth = 15.5703125
$ws2.Columns.Item(2).ColumnWidth = 16.5703125
$ws2.Columns.Item(3).ColumnWidth = 16.140625

# Restore the saved selection on Sheet1, then make Sheet2 the active/
# displayed tab with its own saved selection (so tabSelected + activeTab
# land on Sheet2, as in the source file).
[void]$ws1.Range("C21").Select()
[void]$ws2.Range("F20").Select()
